# cases.xlsx re-scrape: scraper.py was re-run, which picked up a new
# "North America" row for the continent table (inserted at the top, pushing
# the existing five rows down by one) and refreshed every case count to the
# newly scraped figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new first row; everything below shifts down one row.
$ws.Rows.Item(1).Insert()

$continents = @("North America", "Asia", "South America", "Europe", "Africa", "Oceania")
$counts     = @("56,827,972", "80,098,585", "38,549,148", "66,352,163", "8,608,324", "325,319")

# The counts are comma-grouped text (not numbers) in the source data. Excel's
# COM layer "smart types" a literal like "56,827,972" into the number
# 56827972 unless the cell is already formatted as Text, so mark the column
# as Text first. ClearFormats afterwards drops that temporary formatting
# again so no stray number-format sticks around on the finished cells -
# only the literal string values remain.
$counts_range = $ws.Range("B1:B6")
$counts_range.NumberFormat = "@"

for ($i = 0; $i -lt $continents.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $continents[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

$counts_range.ClearFormats()
